$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.241.72'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '3.425.17'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '413.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.72'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.625'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -3.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.726'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.140'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.68'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000221'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.21'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('D14').Value = '3.981.00'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.51'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.58%  '
$ws.Range('D17').Value = '3.438.17'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.78'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +5.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.07'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').Value = '62.260.91'
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '476.22'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +7.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '91.56'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.26'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.02'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.29'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.69'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +8.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.36'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.76'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.63'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.39%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.84'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.73%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.64'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.167'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('E33').Value = '  -3.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '41.10'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.29'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0488'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.70%  '
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.06'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '147.63'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.23%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.134'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.321'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.32'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.07'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.13%  '
$ws.Range('E45').Value = '  +8.47%  '
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.31'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +16.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.31'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('B49').Value = 'PEPE'
$ws.Range('C49').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D49').Value = '0.0₃0535'
$ws.Range('E49').Value = '  +23.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.17'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '113.39'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.80%  '
